$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2,  709, "-1_türkiye_ülke_millet_ver",        "türkiye, ülke, millet, ver, türk"),
    @(3,  188, "0_cumhurbaşkan_aday_ata_ittifak",    "cumhurbaşkan, aday, ata, ittifak, dr"),
    @(4,  184, "1_yayın_program_canlı_konuk",        "yayın, program, canlı, konuk, sun"),
    @(5,  171, "2_millet_oy_buluş_sandık",           "millet, oy, buluş, sandık, seçim"),
    @(6,  134, "3_konut_deprem_depremzede_temel",    "konut, deprem, depremzede, temel, hastane"),
    @(7,  106, "4_bayram_atatürk_türk_kutlu",        "bayram, atatürk, türk, kutlu, mustafa"),
    @(8,  105, "5_türk_türkiye_milliyetçi_yüzyıl",   "türk, türkiye, milliyetçi, yüzyıl, millet"),
    @(9,  82,  "6_şehit_rahmet_atatürk_an",          "şehit, rahmet, atatürk, an, dönüm"),
    @(10, 62,  "7_teşekkür_başkan_dernek_muhteşem",  "teşekkür, başkan, dernek, muhteşem, misafirperverlikleri"),
    @(11, 52,  "8_iyi_sanatçı_allah_iş",              "iyi, sanatçı, allah, iş, çık"),
    @(12, 48,  "9_gel_buluş_saat_bugün",             "gel, buluş, saat, bugün, bekle"),
    @(14, 37,  "11_esnaf_cadde_ziyaret_genç",        "esnaf, cadde, ziyaret, genç, yoğun"),
    @(15, 30,  "12_basın_açıkla_medya_uygula",       "basın, açıkla, medya, uygula, cemiyet"),
    @(16, 30,  "13_muhteşem_bil_şanlıurfa_van",      "muhteşem, bil, şanlıurfa, van, geleneksel"),
    @(17, 27,  "14_emekli_polis_maaş_hak",           "emekli, polis, maaş, hak, düşük"),
    @(18, 25,  "15_nükleer_santral_enerji_üretim",   "nükleer, santral, enerji, üretim, gaz"),
    @(19, 25,  "16_sığınmacı_kaçak_gönder_kararname","sığınmacı, kaçak, gönder, kararname, kal"),
    @(20, 24,  "17_genç_internet_buluş_telefon",     "genç, internet, buluş, telefon, medya"),
    @(21, 23,  "18_öğretmen_eğitim_okul_engelli",    "öğretmen, eğitim, okul, engelli, ata"),
    @(22, 21,  "19_ırak_terör_örgüt_karşı",          "ırak, terör, örgüt, karşı, kardeş"),
    @(23, 14,  "20_dadaş_öv_güzel_çocuk",            "dadaş, öv, güzel, çocuk, erzurum"),
    @(24, 13,  "21_helikopter_uzay_uçak_kuzey",      "helikopter, uzay, uçak, kuzey, milli"),
    @(25, 12,  "22_acı_çerkes_sürgün_kardeş",        "acı, çerkes, sürgün, kardeş, kayıp")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
